# "Inclusão de mais artigos da RBE"
# Extend the "Estudos Econômicos" tracking sheet: the edition/year grid now
# covers years up to 2020 (columns T:BA) for editions 1-4 (rows 2-5), all
# marked "OK", plus a new annotation row (row 6) with a couple of remarks
# about edition nº 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns T:BA (years) for rows 2:5 (editions 1-4) are all "OK".
$ws.Range("T2:BA5").Value = "OK"

# New row 6 with follow-up remarks tied to specific year columns.
$ws.Range("AT6").Value = "INCLUIR 0 DEPOIS"
$ws.Range("AW6").Value = "Todos os artigos do nº 4 são em inglês"
$ws.Range("BA6").Value = "INCLUIR 0 DEPOIS"

# Scroll the view toward the newly filled-in columns and leave BA2 selected
# (mirrors the author's final cursor position after adding the data).
$excel.ActiveWindow.ScrollColumn = 31
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("BA2").Select()
